# Added ifoCAST full series evaluation:
#  - two new quarter columns (Q8, Q9) with headers in row 1
#  - the corresponding diagonal forecast-error values for rows 19-24

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (J1:K1), matching the bold/centered/bordered style
#     already used by the rest of the header row ---
$ws.Range("J1").Value = "Q8"
$ws.Range("K1").Value = "Q9"

$ws.Range("I1").Copy()
$ws.Range("J1:K1").PasteSpecial(-4122)

# --- New diagonal data values, continuing the existing staircase pattern ---
$ws.Range("H19").Value = 2.521807318894689
$ws.Range("I19").Value = 9.990957649751982
$ws.Range("J19").Value = -8.180304150385487
$ws.Range("K19").Value = -0.04063260340633762

$ws.Range("G20").Value = 2.521807318894706
$ws.Range("H20").Value = 9.990957649752
$ws.Range("I20").Value = -8.180304150385469
$ws.Range("J20").Value = -0.04063260340632052

$ws.Range("F21").Value = 2.571807318894714
$ws.Range("G21").Value = 10.04095764975201
$ws.Range("H21").Value = -8.130304150385461
$ws.Range("I21").Value = 0.009367396593688015

$ws.Range("E22").Value = 2.521807318894706
$ws.Range("F22").Value = 9.990957649752
$ws.Range("G22").Value = -8.180304150385469
$ws.Range("H22").Value = -0.04063260340632052

$ws.Range("D23").Value = 2.471807318894706
$ws.Range("E23").Value = 10.039957649752
$ws.Range("F23").Value = -8.131304150385466
$ws.Range("G23").Value = 0.008367396593683296
$ws.Range("H23").Value = 2.155937868393076
$ws.Range("I23").Value = -1.284092118130254
$ws.Range("J23").Value = -1.345266296544455
$ws.Range("K23").Value = 0.6961896406117992

$ws.Range("C24").Value = 2.461807318894701
$ws.Range("D24").Value = 9.980957649751991
$ws.Range("E24").Value = -8.150304150385457
$ws.Range("F24").Value = -0.0006326034063171004
$ws.Range("G24").Value = 2.206937868393077
$ws.Range("H24").Value = -1.233092118130252
$ws.Range("I24").Value = -1.314266296544464
$ws.Range("J24").Value = 0.6871896406117988
